$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Damian Lillard"
$ws.Range("B2").Value = "PG"
$ws.Range("C2").Value = "Milwaukee Bucks"
$ws.Range("A3").Value = "Cade Cunningham"
$ws.Range("B3").Value = "PG,SG"
$ws.Range("C3").Value = "Detroit Pistons"
$ws.Range("A4").Value = "Derrick White"
$ws.Range("B4").Value = "PG,SG"
$ws.Range("C4").Value = "Boston Celtics"
$ws.Range("A5").Value = "Cameron Johnson"
$ws.Range("B5").Value = "SF,PF"
$ws.Range("C5").Value = "Brooklyn Nets"
$ws.Range("A6").Value = "Julius Randle"
$ws.Range("B6").Value = "PF"
$ws.Range("C6").Value = "Minnesota Timberwolves"
$ws.Range("A7").Value = "Kentavious Caldwell-Pope"
$ws.Range("B7").Value = "SG,SF"
$ws.Range("C7").Value = "Orlando Magic"
$ws.Range("A8").Value = "Naz Reid"
$ws.Range("B8").Value = "PF,C"
$ws.Range("C8").Value = "Minnesota Timberwolves"
$ws.Range("A9").Value = "Bam Adebayo"
$ws.Range("B9").Value = "C"
$ws.Range("C9").Value = "Miami Heat"
$ws.Range("A10").Value = "Anthony Davis"
$ws.Range("B10").Value = "PF,C"
$ws.Range("C10").Value = "Los Angeles Lakers"
$ws.Range("A11").Value = "Nick Richards"
$ws.Range("B11").Value = "C"
$ws.Range("C11").Value = "Charlotte Hornets"
$ws.Range("A12").Value = "Isaiah Hartenstein"
$ws.Range("B12").Value = "C"
$ws.Range("C12").Value = "Oklahoma City Thunder"
$ws.Range("A13").Value = "Quentin Grimes"
$ws.Range("B13").Value = "SG,SF"
$ws.Range("C13").Value = "Dallas Mavericks"
$ws.Range("A14").Value = "Malik Monk"
$ws.Range("B14").Value = "SG,SF"
$ws.Range("C14").Value = "Sacramento Kings"
$ws.Range("A15").Value = "Harrison Barnes"
$ws.Range("B15").Value = "SF,PF"
$ws.Range("C15").Value = "San Antonio Spurs"
$ws.Range("A16").Value = "Brandon Miller"
$ws.Range("B16").Value = "SG,SF"
$ws.Range("C16").Value = "Charlotte Hornets"
$ws.Range("A17").Value = "LaMelo Ball"
$ws.Range("B17").Value = "PG,SG"
$ws.Range("C17").Value = "Charlotte Hornets"
$ws.Range("A18").Value = "Brandon Ingram"
$ws.Range("B18").Value = "SG,SF,PF"
$ws.Range("C18").Value = "New Orleans Pelicans"
$ws.Range("A19").Value = "Malcolm Brogdon"
$ws.Range("B19").Value = "PG,SG"
$ws.Range("C19").Value = "Washington Wizards"
